$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new tasks were added to the plan ("Ms. Pac-Man oyununun araştırılması"
# and "Önceki yarışmaların incelenmesi"), pushing every later task down by
# one or two rows, and most of the duration ("Süre (Gün)") values were
# re-tuned. Rewrite the task rows' contents directly (instead of physically
# inserting rows - which, in this host, strips the row's border formatting)
# so every row keeps the "Task" cell formatting (font/fill/border) it
# already has. Rows 12/13 did not have that formatting yet (they used to be
# mostly empty), so copy it over from row 11 first.

# New text is entered in the same order the author must have typed it so
# the shared-string table comes out in the same order as the target file:
# the "Önceki yarışmaların..." task (ends up on row 5) before the
# "Ms. Pac-Man oyununun..." task (ends up on row 3).
$ws.Range("A5").Value = "Önceki yarışmaların incelenmesi"
$ws.Range("A3").Value = "Ms. Pac-Man oyununun araştırılması"

# Rows 12 and 13 used to be mostly empty (just a couple of stray formatted
# cells) - copy the "Task" row formatting from row 11 onto them before
# filling in their real content.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "Proje konusunun kesinleştirilmesi"
$ws.Range("B2").Value = 43724
$ws.Range("C2").Value = 7

$ws.Range("C3").Value = 28

$ws.Range("A4").Value = "Yarışma tarafından sağlanan paketlerin incelenmesi"
$ws.Range("C4").Value = 42

$ws.Range("C5").Value = 21

$ws.Range("A6").Value = "Akıllı algoritmaların araştırılması"
$ws.Range("C6").Value = 21

$ws.Range("A7").Value = "İhtiyaç analizi"
$ws.Range("C7").Value = 21

$ws.Range("A8").Value = "`"Use Case`" diagramlarının oluşturulması"
$ws.Range("C8").Value = 21

$ws.Range("A9").Value = "Ms. Pac-Man ve hayaletler için kullanılacak stratejinin belirlenmesi"
$ws.Range("C9").Value = 14

$ws.Range("A10").Value = "Hayaletler için kullanılacak stratejilerin belirlenmesi"
$ws.Range("C10").Value = 35

$ws.Range("A11").Value = "Yazılımın geliştirilmesi"
$ws.Range("C11").Value = 56

$ws.Range("A12").Value = "Test aşaması ve sonuçların değerlendirilmesi"
$ws.Range("C12").Value = 14

$ws.Range("A13").Value = "Dökümanların hazırlanması"
$ws.Range("C13").Value = 14

# Start-date formulas: B3 on its own, then B4:B13 filled as one shared
# formula (same pattern the original file already used for B3:B11).
$ws.Range("B3").Formula = "=B2 + C2"
$ws.Range("B4:B13").Formula = "=B3 + C3"

# A stray formatted-but-empty cell used to trail row 4 (G4); it now trails
# row 6 instead, two rows further down - same for the other stray cells
# further down the sheet (old E12/A13/B15 -> new E14/A15/B17).
$ws.Range("G4").Clear()
$ws.Range("G6").NumberFormat = "0.00"

$ws.Range("E12").Clear()
$ws.Range("E14").NumberFormat = "0.00"

$ws.Range("A15").NumberFormat = "0.00"

$ws.Range("B15").Clear()
$ws.Range("B17").Value = 43724
$ws.Range("B17").NumberFormat = "0.00"

$ws.Range("G9").Select()
